$d = $word.ActiveDocument
$d.Content.Find.Execute("ni modificaciones al", $true, $false, $false, $false, $false,
                         $true, 1, $false, "en el", 2)
